$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.477.22"
$ws.Range("E2").Value = "  -0.43%  "

$ws.Range("D3").Value = "1.884.93"
$ws.Range("E3").Value = "  -0.58%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'245.93"
$ws.Range("E5").Value = "  -0.69%  "

$ws.Range("D6").Value = "'0.688"
$ws.Range("E6").Value = "  -0.67%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "'42.80"
$ws.Range("E8").Value = "  -3.05%  "

$ws.Range("D9").Value = "'56.50"
$ws.Range("E9").Value = "  +8.49%  "

$ws.Range("E10").Value = "  +0.39%  "

$ws.Range("E11").Value = "  +1.08%  "

$ws.Range("D12").Value = "'0.0983"
$ws.Range("E12").Value = "  +1.43%  "

$ws.Range("D13").Value = "'14.69"
$ws.Range("E13").Value = "  +11.59%  "

$ws.Range("D14").Value = "'0.790"
$ws.Range("E14").Value = "  +7.63%  "

$ws.Range("D15").Value = "2.162.68"
$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("D16").Value = "'5.01"
$ws.Range("E16").Value = "  +0.61%  "

$ws.Range("D17").Value = "1.885.18"
$ws.Range("E17").Value = "  -1.05%  "

$ws.Range("D18").Value = "35.473.96"
$ws.Range("E18").Value = "  -0.39%  "

$ws.Range("D19").Value = "'73.28"
$ws.Range("E19").Value = "  -0.83%  "

$ws.Range("D20").Value = "0.0₃0828"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").Value = "'245.08"
$ws.Range("E21").Value = "  -0.99%  "

$ws.Range("D22").Value = "'13.00"
$ws.Range("E22").Value = "  +0.66%  "

$ws.Range("D23").Value = "'5.16"
$ws.Range("E23").Value = "  +3.94%  "

$ws.Range("D24").Value = "'2.65"
$ws.Range("E24").Value = "  +2.96%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("E26").Value = "  -2.24%  "

$ws.Range("D27").Value = "'165.59"
$ws.Range("E27").Value = "  -0.59%  "

$ws.Range("D28").Value = "'8.61"
$ws.Range("E28").Value = "  +0.79%  "

$ws.Range("D29").Value = "'18.31"
$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").Value = "'0.127"
$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").Value = "'4.40"
$ws.Range("E31").Value = "  +3.19%  "

$ws.Range("D32").Value = "'0.0605"
$ws.Range("E32").Value = "  +3.84%  "

$ws.Range("D33").Value = "'4.26"
$ws.Range("E33").Value = "  +0.38%  "

$ws.Range("D34").Value = "'1.86"
$ws.Range("E34").Value = "  +18.88%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("E36").Value = "  -16.19%  "

$ws.Range("D37").Value = "'0.852"
$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("D38").Value = "'0.0748"
$ws.Range("E38").Value = "  +7.60%  "

$ws.Range("E39").Value = "  -4.30%  "

$ws.Range("E40").Value = "  +6.48%  "

$ws.Range("D41").Value = "'98.63"
$ws.Range("E41").Value = "  +0.70%  "

$ws.Range("D42").Value = "'16.96"
$ws.Range("E42").Value = "  -1.62%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.09"
$ws.Range("E43").Value = "  -0.79%  "

$ws.Range("B44").Value = "Gas"
$ws.Range("C44").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D44").Value = "'14.31"
$ws.Range("E44").Value = "  +18.33%  "

$ws.Range("D45").Value = "1.308.88"
$ws.Range("E45").Value = "  +0.91%  "

$ws.Range("E46").Value = "  -1.97%  "

$ws.Range("E47").Value = "  -1.04%  "

$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("E49").Value = "  -0.83%  "

$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("D51").Value = "'42.52"
$ws.Range("E51").Value = "  -2.15%  "
